$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-29 down to 7-30.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the latest weekly price entry.
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44685
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100112045
$ws.Range("G6").Value = "Zapallo"
$ws.Range("H6").Value = "Camote"
$ws.Range("I6").Value = "1a (cosecha)"
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 680
$ws.Range("L6").Value = 700
$ws.Range("M6").Value = 690
$ws.Range("N6").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O6").Value = "Región de O'Higgins"
$ws.Range("P6").Value = 690
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Hortaliza"
